# Add an "ercc spike-in" column (L) to the processing-metadata template so
# users can flag RNA-seq samples that were spiked with ERCC controls, and
# move the active selection to the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in column L, row 1.
$ws.Range("L1").Value = "with ercc spike-in"

# Flag the two RNA-seq rows: Sample_3 = yes, Sample_4 = no.
$ws.Range("L4").Value = "yes"
$ws.Range("L5").Value = "no"

# Match the author's final selection/active cell.
$null = $ws.Range("L13").Select()
